$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# Row 40
$ws.Range("H40").Value = 1620.8572
$ws.Range("I40").Value = 1389.6
$ws.Range("J40").Value = 2199
$ws.Range("K40").Value = 1389.6
$ws.Range("L40").Value = 2199
$ws.Range("M40").Value = -1214.6
$ws.Range("N40").Value = -2549

# Row 98
$ws.Range("H98").Value = 461.72726
$ws.Range("I98").Value = 358.1
$ws.Range("K98").Value = 358.1
$ws.Range("M98").Value = 1139.9

# Row 107
$ws.Range("H107").Value = 566.4666999999999
$ws.Range("I107").Value = 577.53845
$ws.Range("K107").Value = 577.53845
$ws.Range("M107").Value = 1342.46155

# Row 112
$ws.Range("H112").Value = 0
$ws.Range("I112").Value = 0
$ws.Range("K112").Value = 0
$ws.Range("M112").ClearContents()

# Row 116
$ws.Range("H116").Value = 2521
$ws.Range("I116").Value = 2102.5
$ws.Range("J116").Value = 2800
$ws.Range("K116").Value = 2102.5
$ws.Range("L116").Value = 2800
$ws.Range("M116").Value = 1339.5
$ws.Range("N116").Value = -9684

# Row 122
$ws.Range("H122").Value = 461.72726
$ws.Range("I122").Value = 358.1
$ws.Range("K122").Value = 1074.3
$ws.Range("M122").Value = 1375.7

# Row 132
$ws.Range("H132").Value = 933
$ws.Range("I132").Value = 1049
$ws.Range("J132").Value = 585
$ws.Range("K132").Value = 3147
$ws.Range("L132").Value = 1755
$ws.Range("M132").Value = -617
$ws.Range("N132").Value = -6815

# Row 133
$ws.Range("H133").Value = 0
$ws.Range("J133").Value = 0
$ws.Range("L133").Value = 0
$ws.Range("N133").ClearContents()

$ws = $wb.Worksheets.Item("ARM")
# Row 2
$ws.Range("H2").Value = 1445.6923
$ws.Range("J2").Value = 1522.6666
$ws.Range("L2").Value = 1522.6666
$ws.Range("N2").Value = -1748.6666

# Row 32
$ws.Range("H32").Value = 1881.5
$ws.Range("I32").Value = 1881.5
$ws.Range("K32").Value = 1881.5
$ws.Range("M32").Value = -1594.5

# Row 45
$ws.Range("I45").Value = 2398.8
$ws.Range("K45").Value = 2398.8
$ws.Range("M45").Value = -2021.8

# Row 61
$ws.Range("H61").Value = 2402
$ws.Range("I61").Value = 2312
$ws.Range("J61").Value = 2627
$ws.Range("K61").Value = 2312
$ws.Range("L61").Value = 2627
$ws.Range("M61").Value = -2100
$ws.Range("N61").Value = -3051

# Row 74
$ws.Range("H74").Value = 2784.6
$ws.Range("I74").Value = 2749.5557
$ws.Range("K74").Value = 2749.5557
$ws.Range("M74").Value = -1875.5557

# Row 77
$ws.Range("H77").Value = 2784.6
$ws.Range("I77").Value = 2749.5557
$ws.Range("K77").Value = 13747.7785
$ws.Range("M77").Value = -9379.7785

# Row 116
$ws.Range("H116").Value = 1445.6923
$ws.Range("J116").Value = 1522.6666
$ws.Range("L116").Value = 1522.6666
$ws.Range("N116").Value = -6110.6666

# Row 136
$ws.Range("H136").Value = 2402
$ws.Range("I136").Value = 2312
$ws.Range("J136").Value = 2627
$ws.Range("K136").Value = 6936
$ws.Range("L136").Value = 7881
$ws.Range("M136").Value = -4386
$ws.Range("N136").Value = -12981

$ws = $wb.Worksheets.Item("BSM")
# Row 3
$ws.Range("H3").Value = 1445.6923
$ws.Range("J3").Value = 1522.6666
$ws.Range("L3").Value = 1522.6666
$ws.Range("N3").Value = -1750.6666

# Row 22
$ws.Range("H22").Value = 797.4
$ws.Range("I22").Value = 797.4
$ws.Range("J22").Value = 0
$ws.Range("K22").Value = 797.4
$ws.Range("L22").Value = 0
$ws.Range("M22").Value = -624.4
$ws.Range("N22").ClearContents()

$ws = $wb.Worksheets.Item("CRP")
# Row 4
$ws.Range("H4").Value = 127528
$ws.Range("I4").Value = 111
$ws.Range("J4").Value = 170000.33
$ws.Range("K4").Value = 111
$ws.Range("L4").Value = 170000.33
$ws.Range("M4").Value = 1
$ws.Range("N4").Value = -170224.33

# Row 7
$ws.Range("H7").Value = 72.94118
$ws.Range("I7").Value = 72.76922999999999
$ws.Range("J7").Value = 73.5
$ws.Range("K7").Value = 72.76922999999999
$ws.Range("L7").Value = 73.5
$ws.Range("M7").Value = 40.23077000000001
$ws.Range("N7").Value = -299.5

# Row 16
$ws.Range("H16").Value = 839.8
$ws.Range("I16").Value = 554.875
$ws.Range("J16").Value = 1979.5
$ws.Range("K16").Value = 554.875
$ws.Range("L16").Value = 1979.5
$ws.Range("M16").Value = -267.875
$ws.Range("N16").Value = -2553.5

# Row 22
$ws.Range("H22").Value = 1050.125
$ws.Range("I22").Value = 999.8333
$ws.Range("J22").Value = 1201
$ws.Range("K22").Value = 999.8333
$ws.Range("L22").Value = 1201
$ws.Range("M22").Value = -649.8333
$ws.Range("N22").Value = -1901

# Row 33
$ws.Range("H33").Value = 54566.5
$ws.Range("I33").Value = 51850
$ws.Range("J33").Value = 59999.5
$ws.Range("K33").Value = 51850
$ws.Range("L33").Value = 59999.5
$ws.Range("M33").Value = -51471
$ws.Range("N33").Value = -60757.5

# Row 36
$ws.Range("H36").Value = 22105.875
$ws.Range("I36").Value = 11750
$ws.Range("J36").Value = 32461.75
$ws.Range("K36").Value = 11750
$ws.Range("L36").Value = 32461.75
$ws.Range("M36").Value = -11362
$ws.Range("N36").Value = -33237.75

# Row 39
$ws.Range("H39").Value = 18772.555
$ws.Range("I39").Value = 6158.8335
$ws.Range("J39").Value = 44000
$ws.Range("K39").Value = 6158.8335
$ws.Range("L39").Value = 44000
$ws.Range("M39").Value = -5767.8335
$ws.Range("N39").Value = -44782

# Row 40
$ws.Range("H40").Value = 22105.875
$ws.Range("I40").Value = 11750
$ws.Range("J40").Value = 32461.75
$ws.Range("K40").Value = 11750
$ws.Range("L40").Value = 32461.75
$ws.Range("M40").Value = -11590
$ws.Range("N40").Value = -32781.75

# Row 42
$ws.Range("H42").Value = 32354
$ws.Range("I42").Value = 0
$ws.Range("J42").Value = 32354
$ws.Range("K42").Value = 0
$ws.Range("L42").Value = 32354
$ws.Range("N42").Value = -33540
$ws.Range("M42").ClearContents()

# Row 49
$ws.Range("H49").Value = 18772.555
$ws.Range("I49").Value = 6158.8335
$ws.Range("J49").Value = 44000
$ws.Range("K49").Value = 6158.8335
$ws.Range("L49").Value = 44000
$ws.Range("M49").Value = -5976.8335
$ws.Range("N49").Value = -44364

# Row 74
$ws.Range("H74").Value = 51711
$ws.Range("J74").Value = 51711
$ws.Range("L74").Value = 51711
$ws.Range("N74").Value = -53459

# Row 77
$ws.Range("H77").Value = 51711
$ws.Range("J77").Value = 51711
$ws.Range("L77").Value = 155133
$ws.Range("N77").Value = -163869

# Row 113
$ws.Range("H113").Value = 839.8
$ws.Range("I113").Value = 554.875
$ws.Range("J113").Value = 1979.5
$ws.Range("K113").Value = 554.875
$ws.Range("L113").Value = 1979.5
$ws.Range("M113").Value = 1615.125
$ws.Range("N113").Value = -6319.5

# Row 122
$ws.Range("H122").Value = 831.3333
$ws.Range("I122").Value = 822
$ws.Range("K122").Value = 2466
$ws.Range("M122").Value = -16

# Row 134
$ws.Range("H134").Value = 1326.6666
$ws.Range("I134").Value = 1294.2858
$ws.Range("K134").Value = 3882.8574
$ws.Range("M134").Value = -1347.8574

$ws = $wb.Worksheets.Item("GSM")
# Row 9
$ws.Range("H9").Value = 282.4
$ws.Range("I9").Value = 204.33333
$ws.Range("J9").Value = 399.5
$ws.Range("K9").Value = 204.33333
$ws.Range("L9").Value = 399.5
$ws.Range("M9").Value = -34.33332999999999
$ws.Range("N9").Value = -739.5

# Row 35
$ws.Range("H35").Value = 15007.5
$ws.Range("I35").Value = 15007.5
$ws.Range("K35").Value = 15007.5
$ws.Range("M35").Value = -14709.5

# Row 102
$ws.Range("H102").Value = 941.5
$ws.Range("I102").Value = 848.4167
$ws.Range("K102").Value = 848.4167
$ws.Range("M102").Value = 773.5833

# Row 113
$ws.Range("H113").Value = 743
$ws.Range("I113").Value = 739.5
$ws.Range("J113").Value = 750
$ws.Range("K113").Value = 739.5
$ws.Range("L113").Value = 750
$ws.Range("M113").Value = 1430.5
$ws.Range("N113").Value = -5090

# Row 122
$ws.Range("H122").Value = 4235.5625
$ws.Range("I122").Value = 2801.9092
$ws.Range("J122").Value = 7389.6
$ws.Range("K122").Value = 8405.7276
$ws.Range("L122").Value = 22168.8
$ws.Range("M122").Value = -5955.7276
$ws.Range("N122").Value = -27068.8

# Row 140
$ws.Range("H140").Value = 0
$ws.Range("I140").Value = 0
$ws.Range("K140").Value = 0
$ws.Range("M140").ClearContents()

$ws = $wb.Worksheets.Item("LTW")
# Row 46
$ws.Range("H46").Value = 594
$ws.Range("I46").Value = 598.087
$ws.Range("J46").Value = 500
$ws.Range("K46").Value = 598.087
$ws.Range("L46").Value = 500
$ws.Range("M46").Value = -410.087
$ws.Range("N46").Value = -876

# Row 93
$ws.Range("H93").Value = 3270
$ws.Range("I93").Value = 2950
$ws.Range("J93").Value = 3750
$ws.Range("K93").Value = 2950
$ws.Range("L93").Value = 3750
$ws.Range("M93").Value = -1702
$ws.Range("N93").Value = -6246

# Row 122
$ws.Range("H122").Value = 3785.6428
$ws.Range("I122").Value = 3533.1667
$ws.Range("J122").Value = 3975
$ws.Range("K122").Value = 10599.5001
$ws.Range("L122").Value = 11925
$ws.Range("M122").Value = -8149.500100000001
$ws.Range("N122").Value = -16825

# Row 132
$ws.Range("H132").Value = 3099.75
$ws.Range("I132").Value = 3099.75
$ws.Range("K132").Value = 9299.25
$ws.Range("M132").Value = -6769.25

$ws = $wb.Worksheets.Item("WVR")
# Row 122
$ws.Range("H122").Value = 999
$ws.Range("I122").Value = 999
$ws.Range("K122").Value = 2997
$ws.Range("M122").Value = -547

# Row 126
$ws.Range("H126").Value = 4483.25
$ws.Range("I126").Value = 2999.875
$ws.Range("J126").Value = 7450
$ws.Range("K126").Value = 8999.625
$ws.Range("L126").Value = 22350
$ws.Range("M126").Value = -6529.625
$ws.Range("N126").Value = -27290

# Row 132
$ws.Range("H132").Value = 1020.3333
$ws.Range("I132").Value = 1022.875
$ws.Range("J132").Value = 1000
$ws.Range("K132").Value = 3068.625
$ws.Range("L132").Value = 1000
$ws.Range("M132").Value = -538.625
$ws.Range("N132").Value = -8060
